# Daily attendance processing - 2026-01-18 15:34:38
# Swap the "Recorded By" value ordering for rows that contain both
# "System" and "dnasr281@gmail.com" - from "System, dnasr281@gmail.com"
# to "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$lastRow = $ws.Cells.SpecialCells(11).Row  # xlCellTypeLastCell = 11

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
